# Refresh "想去人数" (interest-count) figures in column F across all four
# sheets of the 北京-漫展信息 workbook, matching the regenerated gh-pages
# data snapshot (commit 456a3b4). Only column F values change; everything
# else (labels, links, prices, etc.) stays untouched.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 3530
$ws.Range("F5").Value  = 3530
$ws.Range("F6").Value  = 252
$ws.Range("F7").Value  = 5056
$ws.Range("F8").Value  = 510
$ws.Range("F9").Value  = 345
$ws.Range("F11").Value = 681
$ws.Range("F13").Value = 77
$ws.Range("F16").Value = 308
$ws.Range("F17").Value = 32
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 359
$ws.Range("F22").Value = 4886
$ws.Range("F26").Value = 6003
$ws.Range("F29").Value = 3217
$ws.Range("F30").Value = 327
$ws.Range("F31").Value = 701
$ws.Range("F32").Value = 4440
$ws.Range("F34").Value = 119
$ws.Range("F35").Value = 139
$ws.Range("F36").Value = 991
$ws.Range("F37").Value = 81
$ws.Range("F40").Value = 855
$ws.Range("F41").Value = 960
$ws.Range("F42").Value = 16

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 45

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 225
$ws.Range("F3").Value = 1112

# ---- 全部类型 (All types, union of the above sheets) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 225
$ws.Range("F4").Value  = 1112
$ws.Range("F8").Value  = 3530
$ws.Range("F9").Value  = 3530
$ws.Range("F10").Value = 252
$ws.Range("F11").Value = 5056
$ws.Range("F12").Value = 510
$ws.Range("F13").Value = 345
$ws.Range("F15").Value = 681
$ws.Range("F16").Value = 77
$ws.Range("F19").Value = 308
$ws.Range("F20").Value = 32
$ws.Range("F21").Value = 45
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 359
$ws.Range("F26").Value = 4886
$ws.Range("F30").Value = 6003
$ws.Range("F33").Value = 3217
$ws.Range("F34").Value = 327
$ws.Range("F35").Value = 701
$ws.Range("F36").Value = 4440
$ws.Range("F39").Value = 119
$ws.Range("F40").Value = 139
$ws.Range("F41").Value = 991
$ws.Range("F42").Value = 81
$ws.Range("F45").Value = 855
$ws.Range("F46").Value = 960
$ws.Range("F48").Value = 16
